$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing "micro-rna" filename in column A for the MicroRNAs
# in Type 1 Diabetes row (row 10).
$ws.Range("A10").Value = "micro-rna"

# Move the selection to the newly filled cell, with the view scrolled
# back to the top-left (A1).
$ws.Range("A1").Select()
$ws.Range("A10").Select()
